# Fruta / hortaliza, semanal
# Weekly update: insert a new price record as row 317, shifting the
# existing rows 317-402 down to 318-403 (dimension grows to A1:R403).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 317; this pushes the old rows 317..402
# down to 318..403 and inherits the formatting (e.g. the date style on
# column D) from the surrounding rows.
$ws.Rows.Item(317).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A317").Value = 5
$ws.Range("B317").Value = "Macroferia Regional de Talca"
$ws.Range("C317").Value = "Maule"
$ws.Range("D317").Value = 44551
$ws.Range("E317").Value = 7
$ws.Range("F317").Value = 100112002
$ws.Range("G317").Value = "Pimiento"
$ws.Range("H317").Value = "Cuatro cascos verde"
$ws.Range("I317").Value = "Primera"
$ws.Range("J317").Value = 200
$ws.Range("K317").Value = 12000
$ws.Range("L317").Value = 12000
$ws.Range("M317").Value = 12000
$ws.Range("N317").Value = "`$/caja 15 kilos"
$ws.Range("O317").Value = "Región del Maule"
$ws.Range("P317").Value = 800
$ws.Range("Q317").Value = 15
$ws.Range("R317").Value = "Hortaliza"
